# Update "Pais" (countries) COVID data sheet and a few provincias-style
# reorderings in the country list, matching a newer data pull (17:05 -> 18:05).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / last-updated timestamp
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 18:05"
$ws.Range("B4").Value = 1799225
$ws.Range("C4").Value = 5695
$ws.Range("D4").Value = 519736
$ws.Range("E4").Value = 1174809
$ws.Range("G4").Value = 138
$ws.Range("H4").Value = 104680
$ws.Range("B5").Value = 469510
$ws.Range("C5").Value = 1172
$ws.Range("E5").Value = 248314
$ws.Range("G5").Value = 71
$ws.Range("H5").Value = 28015
$ws.Range("B8").Value = 272826
$ws.Range("C8").Value = 1604
$ws.Range("G8").Value = 215
$ws.Range("H8").Value = 38376
$ws.Range("B11").Value = 183139
$ws.Range("C11").Value = 120
$ws.Range("E11").Value = 9641
$ws.Range("B12").Value = 180939
$ws.Range("C12").Value = 7448
$ws.Range("D12").Value = 84792
$ws.Range("E12").Value = 91003
$ws.Range("G12").Value = 164
$ws.Range("H12").Value = 5144
$ws.Range("B16").Value = 94858
$ws.Range("C16").Value = 4220
$ws.Range("D16").Value = 40431
$ws.Range("E16").Value = 53430
$ws.Range("G16").Value = 53
$ws.Range("H16").Value = 997
$ws.Range("B38").Value = 23571
$ws.Range("C38").Value = 416
$ws.Range("E38").Value = 11494
$ws.Range("G38").Value = 10
$ws.Range("H38").Value = 1061
$ws.Range("D55").Value = 5220
$ws.Range("E55").Value = 5124

# Argelia / Chequia swap positions (rows 57-58) with refreshed data
$ws.Range("A57").Value = "Argelia"
$ws.Range("B57").Value = 9267
$ws.Range("C57").Value = 133
$ws.Range("D57").Value = 5549
$ws.Range("E57").Value = 3072
$ws.Range("G57").Value = 8
$ws.Range("H57").Value = 646
$ws.Range("A58").Value = "Chequia"
$ws.Range("B58").Value = 9226
$ws.Range("C58").Value = 30
$ws.Range("D58").Value = 6532
$ws.Range("E58").Value = 2375
$ws.Range("H58").Value = 319
$ws.Range("B69").Value = 5659
$ws.Range("C69").Value = 223
$ws.Range("D69").Value = 3441
$ws.Range("E69").Value = 2033
$ws.Range("G69").Value = 8
$ws.Range("H69").Value = 185
$ws.Range("B74").Value = 4016
$ws.Range("C74").Value = 4
$ws.Range("E74").Value = 91
$ws.Range("B83").Value = 2915
$ws.Range("C83").Value = 6
$ws.Range("E83").Value = 1366
$ws.Range("D92").Value = 464
$ws.Range("E92").Value = 1361
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 63
$ws.Range("D126").Value = 692
$ws.Range("E126").Value = 21
$ws.Range("B131").Value = 734
$ws.Range("C131").Value = 4
$ws.Range("E131").Value = 218

# Comoras moves up (rows 173-177 shift down one position) with refreshed data
$ws.Range("A173").Value = "Comoras"
$ws.Range("B173").Value = 106
$ws.Range("C173").Value = 19
$ws.Range("D173").Value = 26
$ws.Range("E173").Value = 78
$ws.Range("H173").Value = 2
$ws.Range("A174").Value = "Bahamas"
$ws.Range("B174").Value = 102
$ws.Range("D174").Value = 48
$ws.Range("E174").Value = 43
$ws.Range("H174").Value = 11
$ws.Range("A175").Value = "Aruba"
$ws.Range("B175").Value = 101
$ws.Range("D175").Value = 98
$ws.Range("E175").Value = 0
$ws.Range("H175").Value = 3
$ws.Range("A176").Value = "Monaco"
$ws.Range("B176").Value = 98
$ws.Range("D176").Value = 90
$ws.Range("E176").Value = 4
$ws.Range("H176").Value = 4
$ws.Range("A177").Value = "Barbados"
$ws.Range("B177").Value = 92
$ws.Range("D177").Value = 76
$ws.Range("E177").Value = 9
$ws.Range("H177").Value = 7

# Belice / Santa Lucia swap positions (rows 200-201)
$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

# San Bartolome / Bonaire, San Eustaquio y Saba swap positions (rows 215-216)
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "San Bartolome"
